# Generate Report for Handoff
#
# Updates the localization-status report to reflect that the handoff
# report has just been (re-)generated:
#   - Status values move from "Handed back: in sync with en-US" to
#     "Ready for handoff" on the Overview sheet and on each language
#     sheet's Status column.
#   - The "Latest Handoff Datetime" / "Latest HO Xliff Generate Date"
#     timestamps are refreshed to the new generation time.
#   - The Status columns, which previously needed to fit the long
#     "Handed back: in sync with en-US" text, are narrowed to fit the
#     shorter "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status text + timestamp updates ---------------------------------

# Overview sheet: per-language status (E2, F2) and generation date (G2)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-20 11:03:28"

# zh-cn sheet: Status (C2) and Latest Handoff Datetime (H2)
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-20 11:03:24"

# de-de sheet: Status (C2) and Latest Handoff Datetime (H2)
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-20 11:03:28"

# --- Column width updates (status columns narrowed) -------------------

# Overview: zh-cn (E) and de-de (F) status columns
$wsOverview.Range("E1:F1").ColumnWidth = 16.33

# zh-cn / de-de: Status column (C)
$wsZhCn.Range("C1").ColumnWidth = 16.33
$wsDeDe.Range("C1").ColumnWidth = 16.33
